$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$m = $s.Master
Write-Output $m
$cs = $m.Theme
Write-Output $cs
